$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.104631
$ws.Range("H2").Value = 3.313893
$ws.Range("I2").Value = 0.8734451962653081
$ws.Range("J2").Value = 0.8734451962653083
$ws.Range("M2").Value = 6.436245333333333
$ws.Range("N2").Value = 19.308736
$ws.Range("O2").Value = 0.2367562936388591
$ws.Range("P2").Value = 0.2367562936388591
$ws.Range("Q2").Value = 7.109676118805334
$ws.Range("R2").Value = 63.987085069248
$ws.Range("S2").Value = 0.2067936473644402
$ws.Range("T2").Value = 0.2067936473644403

# Row 3 (ECs -> FAPs)
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.104631
$ws.Range("H3").Value = 3.313893
$ws.Range("I3").Value = 0.8734451962653081
$ws.Range("J3").Value = 0.8734451962653083
$ws.Range("O3").Value = 0.5508630013028089
$ws.Range("P3").Value = 0.550863001302809
$ws.Range("Q3").Value = 16.54214747537
$ws.Range("R3").Value = 148.87932727833
$ws.Range("S3").Value = 0.4811486422882286
$ws.Range("T3").Value = 0.4811486422882287

# Row 4 (ECs -> MuSCs)
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.104631
$ws.Range("H4").Value = 3.313893
$ws.Range("I4").Value = 0.8734451962653081
$ws.Range("J4").Value = 0.8734451962653083
$ws.Range("O4").Value = 0.212380705058332
$ws.Range("P4").Value = 0.212380705058332
$ws.Range("Q4").Value = 6.377689072762334
$ws.Range("R4").Value = 57.399201654861
$ws.Range("S4").Value = 0.1855029066126393
$ws.Range("T4").Value = 0.1855029066126394

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.1265548037346918
$ws.Range("J5").Value = 0.1265548037346918
$ws.Range("M5").Value = 6.436245333333333
$ws.Range("N5").Value = 19.308736
$ws.Range("O5").Value = 0.2367562936388591
$ws.Range("P5").Value = 0.2367562936388591
$ws.Range("Q5").Value = 1.030131792675556
$ws.Range("R5").Value = 9.271186134080001
$ws.Range("S5").Value = 0.02996264627441887
$ws.Range("T5").Value = 0.02996264627441888

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.1265548037346918
$ws.Range("J6").Value = 0.1265548037346918
$ws.Range("O6").Value = 0.5508630013028089
$ws.Range("P6").Value = 0.550863001302809
$ws.Range("S6").Value = 0.06971435901458024
$ws.Range("T6").Value = 0.06971435901458026

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.1265548037346918
$ws.Range("J7").Value = 0.1265548037346918
$ws.Range("O7").Value = 0.212380705058332
$ws.Range("P7").Value = 0.212380705058332
$ws.Range("S7").Value = 0.02687779844569267
$ws.Range("T7").Value = 0.02687779844569268
